$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update report generation timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:02 AM"

# Update Total Billed Amount
$ws.Range("C8").Value = 478.55

# Clear Scope ID # value (was "#NO MATCH", now blank)
$ws.Range("G10").Value = ""

# Update line item pricing
$ws.Range("H16").Value = 478.55

# Update total pricing
$ws.Range("H17").Value = 478.55
